# Product Backlog update following Sprint 1 demo:
# re-prioritize a handful of backlog items, then re-sort the backlog
# by Priority (column A) ascending, exactly as "Data > Sort" would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-prioritize specific backlog items (pre-sort positions) ---
$ws.Range("A2").Value  = 0.5   # Create a Trip (name and start/end time)
$ws.Range("A3").Value  = 0.5   # Add Waypoints (place + time) to my trip
$ws.Range("A7").Value  = 2     # Add Lodging to a trip
$ws.Range("A8").Value  = 2     # Remove Lodging from a trip
$ws.Range("A9").Value  = 0.5   # View an overview of my trip
$ws.Range("A14").Value = 0.5   # Add notes to an item in my Trip

# --- Re-sort the backlog by Priority ascending (Data > Sort), same
#     range the sheet's saved sort state already used ---
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A2:A29"))
$so.SetRange($ws.Range("A2:D29"))
$so.Header = -4142   # xlNo
$so.Apply()

# --- View state: move the selection and zoom in ---
$ws.Range("A3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
